# Inductor Design_CA.xlsx - "Report is added, transformer parameters are recalculated."
#
# The raw rows 18-25 (measured datasheet samples) and the derived summary
# row (28-29) are cleared back to blank inputs/formulas, which in turn
# drops the now-unused "B_mean(T)" shared string and renumbers the
# sharedStrings table. The sheet's active selection is also moved to
# reflect the newly-cleared report range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the calculated/measured block C18:H25 (formulas + values), leaving
# the cells blank but keeping their existing styles.
$ws.Range("C18:H25").Value = ""

# Clear the "B_mean(T)" / "u_0" labels on row 28 ...
$ws.Range("C28").Value = ""
$ws.Range("E28").Value = ""

# ... and the averaged summary formulas on row 29.
$ws.Range("C29:E29").Value = ""

# Move the active selection to the cleared report range.
$null = $ws.Range("C18:H29").Select()
